$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "No." column before the existing column B -------------
# This shifts the old B:G columns (and all their data/formatting) one
# column to the right, becoming C:H.
$ws.Columns("B:B").Insert()

# --- Fix up the cell that used to hold the "Waiting" status -------------
# Column G (old) / H (new), row 10 used to be "Waiting" (style 4); the
# refreshed tracker marks it "Done" using the same look as the other
# "Done" cells in that row band (style copied from G9, which already
# reads "Done"-style formatting).
$ws.Range("G9").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("H10").Value = "Done"

# --- Header for the new column -------------------------------------------
# Build the header style by reusing the existing bordered look (copied
# from B3, which only carries font/border, no alignment) and then
# layering on the bold font + yellow fill that the other header cells
# use - this reconstructs the header style without forcing the centered
# alignment the other headers have.
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Interior.Color = 65535
$ws.Range("B2").Value = "No."

# --- Row numbers 1-8 for the task rows -----------------------------------
# Re-use the plain bordered look already present on G8 (style with just a
# border, no fill) for every numbered cell.
$ws.Range("G8").Copy()
$numberRange = $ws.Range("B3:B10")
$numberRange.PasteSpecial(-4122)

$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 6
$ws.Range("B9").Value = 7
$ws.Range("B10").Value = 8

# --- Column widths ---------------------------------------------------------
$ws.Columns("B:B").ColumnWidth = 3.25
$ws.Columns("C:C").ColumnWidth = 27.5703125

# --- Row heights for the two wrapped/multi-line task rows ------------------
$ws.Rows("9:9").RowHeight = 90
$ws.Rows("10:10").RowHeight = 90

# --- Selection, matching the saved view in the workbook --------------------
$ws.Range("C9").Select()
